# Updated cryptos list - apply price (D) and volume change (E) edits
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to keep the assigned value as literal text,
    # even when it looks like a number (e.g. "583.68"), then restore
    # the original (default) cell style so no formatting is changed.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "63.699.94"
Set-TextValue $ws.Range("E2") "  -1.87%  "
Set-TextValue $ws.Range("D3") "3.490.08"
Set-TextValue $ws.Range("E3") "  -0.87%  "
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "583.68"
Set-TextValue $ws.Range("E5") "  -2.20%  "
Set-TextValue $ws.Range("D6") "130.94"
Set-TextValue $ws.Range("E6") "  -1.84%  "
Set-TextValue $ws.Range("D7") "3.491.08"
Set-TextValue $ws.Range("E7") "  -0.81%  "
Set-TextValue $ws.Range("E8") "  +0.02%  "
Set-TextValue $ws.Range("E9") "  -1.82%  "
Set-TextValue $ws.Range("E10") "  +0.02%  "
Set-TextValue $ws.Range("D11") "7.14"
Set-TextValue $ws.Range("E11") "  -0.02%  "
Set-TextValue $ws.Range("E12") "  -0.52%  "
Set-TextValue $ws.Range("D13") "4.078.57"
Set-TextValue $ws.Range("E13") "  -0.97%  "
Set-TextValue $ws.Range("D14") "27.31"
Set-TextValue $ws.Range("E14") "  -0.41%  "
Set-TextValue $ws.Range("D15") "0.119"
Set-TextValue $ws.Range("E15") "  +1.49%  "
Set-TextValue $ws.Range("E16") "  -2.23%  "
Set-TextValue $ws.Range("D17") "3.466.05"
Set-TextValue $ws.Range("E17") "  -1.49%  "
Set-TextValue $ws.Range("D18") "63.809.37"
Set-TextValue $ws.Range("E18") "  -1.75%  "
Set-TextValue $ws.Range("D19") "9.86"
Set-TextValue $ws.Range("E19") "  -2.14%  "
Set-TextValue $ws.Range("D20") "14.13"
Set-TextValue $ws.Range("E20") "  -1.24%  "
Set-TextValue $ws.Range("D21") "5.63"
Set-TextValue $ws.Range("E21") "  -0.75%  "
Set-TextValue $ws.Range("D22") "383.08"
Set-TextValue $ws.Range("E22") "  -1.92%  "
Set-TextValue $ws.Range("D24") "3.628.32"
Set-TextValue $ws.Range("E24") "  -0.92%  "
Set-TextValue $ws.Range("D25") "73.28"
Set-TextValue $ws.Range("E25") "  -0.88%  "
Set-TextValue $ws.Range("D27") "0.0000114"
Set-TextValue $ws.Range("E27") "  +2.64%  "
Set-TextValue $ws.Range("D28") "1.58"
Set-TextValue $ws.Range("E28") "  -1.57%  "
Set-TextValue $ws.Range("D29") "7.54"
Set-TextValue $ws.Range("E29") "  -2.30%  "
Set-TextValue $ws.Range("E30") "  +0.13%  "
Set-TextValue $ws.Range("D31") "8.27"
Set-TextValue $ws.Range("E31") "  -0.73%  "
Set-TextValue $ws.Range("E32") "  -2.32%  "
Set-TextValue $ws.Range("D33") "3.496.08"
Set-TextValue $ws.Range("E33") "  -0.81%  "
Set-TextValue $ws.Range("E34") "  -0.01%  "
Set-TextValue $ws.Range("D35") "23.45"
Set-TextValue $ws.Range("E35") "  -3.10%  "
Set-TextValue $ws.Range("E36") "  +0.43%  "
Set-TextValue $ws.Range("D37") "5.33"
Set-TextValue $ws.Range("E37") "  +4.13%  "
Set-TextValue $ws.Range("D38") "6.96"
Set-TextValue $ws.Range("E38") "  +1.82%  "
Set-TextValue $ws.Range("E39") "  -0.75%  "
Set-TextValue $ws.Range("D40") "160.56"
Set-TextValue $ws.Range("E40") "  -4.70%  "
Set-TextValue $ws.Range("D41") "0.0794"
Set-TextValue $ws.Range("E41") "  -2.76%  "
Set-TextValue $ws.Range("D42") "26.64"
Set-TextValue $ws.Range("E42") "  +4.42%  "
Set-TextValue $ws.Range("E43") "  -1.58%  "
Set-TextValue $ws.Range("E44") "  -0.02%  "
Set-TextValue $ws.Range("D45") "41.64"
Set-TextValue $ws.Range("E45") "  -2.50%  "
Set-TextValue $ws.Range("E46") "  -2.65%  "
Set-TextValue $ws.Range("D47") "4.38"
Set-TextValue $ws.Range("E47") "  -0.67%  "
Set-TextValue $ws.Range("E48") "  -1.74%  "
Set-TextValue $ws.Range("E49") "  -0.72%  "
Set-TextValue $ws.Range("D50") "2.421.41"
Set-TextValue $ws.Range("E50") "  +1.61%  "
Set-TextValue $ws.Range("D51") "0.899"
Set-TextValue $ws.Range("E51") "  +1.43%  "
